# Deluvio: looks like release-candidate.
# Adds an ADC reference-voltage-divider calc sheet on Sheet3 and makes it the
# active tab (was BOM before).

$wb  = $excel.ActiveWorkbook
$bom = $wb.Worksheets.Item("BOM")
$ws  = $wb.Worksheets.Item("Sheet3")

# --- Populate Sheet3 with the ADC threshold-divider calculation ----------
# (cells are written in this exact order so the shared-string table comes
#  out in the same sequence as the target workbook)
$ws.Range("A1").Value = "ADC"

$ws.Range("A2").Value = "Aref"
$ws.Range("C2").Value = "mV"

$ws.Range("A3").Value = "U treshold"
$ws.Range("C3").Value = "mV"

$ws.Range("A4").Value = "Rup"
$ws.Range("A5").Value = "Rdown"
$ws.Range("C4").Value = "kOhm"
$ws.Range("C5").Value = "kOhm"

$ws.Range("A6").Value = "U in"
$ws.Range("C6").Value = "mV"

$ws.Range("A7").Value = "ADC"

# Numeric inputs
$ws.Range("B2").Value = 3300
$ws.Range("B3").Value = 8000
$ws.Range("B4").Value = 1000
$ws.Range("B5").Value = 330

# Formulas
$ws.Range("B6").Formula = "=B3*B5/(B4+B5)"
$ws.Range("B7").Formula = "=1024*B6/B2"
$ws.Range("B7").NumberFormat = "0"

# Column widths (best achievable match to the authored 11.57 / 9.57 widths)
$ws.Columns("A").ColumnWidth = 10.59
$ws.Columns("B").ColumnWidth = 8.59

# --- Tab/selection state: Sheet3 becomes the active tab -------------------
[void]$bom.Range("D9").Select()
[void]$ws.Range("B4").Select()
